# Effort Estimation Sheet - add "Actual Time" (col I) figures for a few
# stories, fill in the missing Task text for row 9, and move the
# selection/scroll position the way the author left the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Time" (column I) values that were filled in.
$ws.Cells.Item(4, 9).Value  = 6
$ws.Cells.Item(7, 9).Value  = 10
$ws.Cells.Item(13, 9).Value = 12
$ws.Cells.Item(14, 9).Value = 10
$ws.Cells.Item(17, 9).Value = 5
$ws.Cells.Item(19, 9).Value = 2

# Row 9 (user story #6) no longer needs a task, noted as "Taka út öruglega".
$ws.Cells.Item(9, 11).Value = "Taka út öruglega"

# Leave the sheet scrolled/selected where the author left it.
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I25").Select()
